$ws = $excel.ActiveWorkbook.ActiveSheet

# Pre-format the handful of Price cells whose target text has a trailing
# zero after the decimal point (e.g. "1.000", "6.310") as Text, so Excel
# does not silently renormalize them to a shorter numeric literal.
foreach ($addr in @("D18", "D22", "D25", "D37", "D38", "D42")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "24.586.55"
$ws.Range("E2").Value = "  +1.55%  "
$ws.Range("D3").Value = "1.703.45"
$ws.Range("E3").Value = "  +1.67%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "308.82"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").Value = "0.9993"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "0.3729"
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").Value = "48.84"
$ws.Range("E8").Value = "  +2.67%  "
$ws.Range("D9").Value = "0.3415"
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("D10").Value = "1.176"
$ws.Range("E10").Value = "  -0.61%  "
$ws.Range("D11").Value = "0.07415"
$ws.Range("E11").Value = "  +1.43%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").Value = "20.74"
$ws.Range("E13").Value = "  +1.56%  "
$ws.Range("D14").Value = "6.197"
$ws.Range("E14").Value = "  +1.33%  "
$ws.Range("D15").Value = "6.886"
$ws.Range("E15").Value = "  +1.63%  "
$ws.Range("D16").Value = "1.696.87"
$ws.Range("E16").Value = "  +1.27%  "
$ws.Range("D17").Value = "0.00001113"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").Value = "0.06665"
$ws.Range("E19").Value = "  -1.00%  "
$ws.Range("D21").Value = "16.99"
$ws.Range("E21").Value = "  +2.55%  "
$ws.Range("D22").Value = "6.310"
$ws.Range("E22").Value = "  +2.64%  "
$ws.Range("D23").Value = "13.14"
$ws.Range("E23").Value = "  +9.11%  "
$ws.Range("D24").Value = "24.580.79"
$ws.Range("E24").Value = "  +1.67%  "
$ws.Range("D25").Value = "2.440"
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("D26").Value = "2.752"
$ws.Range("E26").Value = "  +3.28%  "
$ws.Range("D27").Value = "20.09"
$ws.Range("E27").Value = "  +2.79%  "
$ws.Range("D28").Value = "149.22"
$ws.Range("E28").Value = "  -1.56%  "
$ws.Range("D29").Value = "130.71"
$ws.Range("E29").Value = "  +3.09%  "
$ws.Range("D30").Value = "1.884.74"
$ws.Range("E30").Value = "  +1.28%  "
$ws.Range("D31").Value = "1.168"
$ws.Range("E31").Value = "  +17.57%  "
$ws.Range("D32").Value = "6.637"
$ws.Range("E32").Value = "  +2.24%  "
$ws.Range("D33").Value = "4.205"
$ws.Range("E33").Value = "  +2.70%  "
$ws.Range("D34").Value = "0.08733"
$ws.Range("E34").Value = "  +2.51%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").Value = "13.52"
$ws.Range("E35").Value = "  +6.88%  "
$ws.Range("B36").Value = "WEMIXTOKEN"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "1.756"
$ws.Range("E36").Value = "  -0.80%  "
$ws.Range("D37").Value = "5.470"
$ws.Range("E37").Value = "  +1.71%  "
$ws.Range("D38").Value = "0.06470"
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "8.853"
$ws.Range("E39").Value = "  -0.82%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.02357"
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").Value = "0.2177"
$ws.Range("E41").Value = "  +1.49%  "
$ws.Range("D42").Value = "1.270"
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("D43").Value = "0.6371"
$ws.Range("E43").Value = "  +2.69%  "
$ws.Range("D44").Value = "0.9993"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "13.87"
$ws.Range("E45").Value = "  +5.37%  "
$ws.Range("D46").Value = "0.6039"
$ws.Range("E46").Value = "  +1.41%  "
$ws.Range("D47").Value = "3.801"
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("D48").Value = "2.099"
$ws.Range("E48").Value = "  +3.18%  "
$ws.Range("D49").Value = "128.32"
$ws.Range("E49").Value = "  +0.66%  "
$ws.Range("D50").Value = "0.07217"
$ws.Range("E50").Value = "  +0.60%  "
$ws.Range("D51").Value = "78.56"
$ws.Range("E51").Value = "  +2.52%  "
